# Updated cryptos list on Tue May 21 10:43:11 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $text)
    # Price cells are stored as plain text even when they look like numbers
    # (e.g. "620.47"); force text format first so COM doesn't coerce them
    # into numeric cells, then restore the original (default) cell style so
    # no stray formatting is left behind.
    $ws.Range($range).NumberFormat = "@"
    $ws.Range($range).Value = $text
    $ws.Range($range).Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "71.113.30"
$ws.Range("E2").Value = "  +5.93%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.721.07"
$ws.Range("E3").Value = "  +19.88%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
Set-TextValue "D5" "620.47"
$ws.Range("E5").Value = "  +7.96%  "

# Row 6 - Solana
Set-TextValue "D6" "183.22"
$ws.Range("E6").Value = "  +3.10%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.718.26"
$ws.Range("E7").Value = "  +19.85%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.00%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  +5.62%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +8.23%  "

# Row 11 - Toncoin
Set-TextValue "D11" "6.62"
$ws.Range("E11").Value = "  +3.86%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.503"
$ws.Range("E12").Value = "  +7.32%  "

# Row 13 - Avalanche
Set-TextValue "D13" "40.74"
$ws.Range("E13").Value = "  +13.07%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  +6.35%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "4.333.82"
$ws.Range("E15").Value = "  +19.74%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.712.04"
$ws.Range("E16").Value = "  +19.53%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "71.116.74"

# Row 18 - TRON
$ws.Range("E18").Value = "  +1.71%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +7.10%  "

# Row 20 & 21 - swap BitcoinCash / Chainlink
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D20" "16.93"
$ws.Range("E20").Value = "  +1.52%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D21" "519.36"
$ws.Range("E21").Value = "  +4.46%  "

# Row 22 - Uniswap
Set-TextValue "D22" "9.34"
$ws.Range("E22").Value = "  +20.20%  "

# Row 23 - Polygon
Set-TextValue "D23" "0.746"
$ws.Range("E23").Value = "  +8.52%  "

# Row 24 & 25 - swap Litecoin / Fetch.AI
$ws.Range("B24").Value = "Fetch.AI"
$ws.Range("C24").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D24" "2.55"
$ws.Range("E24").Value = "  +12.89%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D25" "88.90"
$ws.Range("E25").Value = "  +6.18%  "

# Row 26 - InternetComputer(DFINITY)
Set-TextValue "D26" "13.59"
$ws.Range("E26").Value = "  +7.85%  "

# Row 27 - RenderToken
Set-TextValue "D27" "11.25"
$ws.Range("E27").Value = "  +11.56%  "

# Row 28 - Dai
$ws.Range("E28").Value = "  -0.01%  "

# Row 29 - ImmutableX
Set-TextValue "D29" "2.56"

# Row 30 - NEARProtocol
Set-TextValue "D30" "8.22"
$ws.Range("E30").Value = "  +3.87%  "

# Row 31 - PancakeSwap
Set-TextValue "D31" "2.90"
$ws.Range("E31").Value = "  +11.88%  "

# Row 32 - EthereumClassic
Set-TextValue "D32" "31.93"
$ws.Range("E32").Value = "  +13.27%  "

# Row 33 - PEPE
$ws.Range("E33").Value = "  +18.18%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  +4.52%  "

# Row 35 - FirstDigitalUSD
Set-TextValue "D35" "0.999"
$ws.Range("E35").Value = "  +0.01%  "

# Row 37 - Mantle
$ws.Range("E37").Value = "  +10.13%  "

# Row 38 & 39 - swap TheGraph / Stacks
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D38" "2.27"
$ws.Range("E38").Value = "  +12.62%  "

$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D39" "0.347"
$ws.Range("E39").Value = "  +10.84%  "

# Row 40 - Kaspa
$ws.Range("E40").Value = "  +9.40%  "

# Row 42 - Bittensor
Set-TextValue "D42" "436.75"
$ws.Range("E42").Value = "  +17.72%  "

# Row 43 - Arweave
Set-TextValue "D43" "45.12"
$ws.Range("E43").Value = "  -5.61%  "

# Row 44 & 45 - swap Cosmos / Maker
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D44" "3.151.57"
$ws.Range("E44").Value = "  +13.06%  "

$ws.Range("B45").Value = "Cosmos"
$ws.Range("C45").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D45" "8.87"
$ws.Range("E45").Value = "  +6.98%  "

# Row 46 - dogwifhat
Set-TextValue "D46" "2.88"
$ws.Range("E46").Value = "  +5.17%  "

# Row 47 - VeChain
$ws.Range("E47").Value = "  +7.33%  "

# Row 48 - InjectiveProtocol
Set-TextValue "D48" "28.34"
$ws.Range("E48").Value = "  +11.47%  "

# Row 49 - Monero
Set-TextValue "D49" "140.56"
$ws.Range("E49").Value = "  +3.82%  "

# Row 51 - ThetaToken
$ws.Range("E51").Value = "  +8.43%  "
